$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Update status text: "Handed back: in sync with en-US" -> "Ready for handoff" ---
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("C2").Value = "Ready for handoff"

# --- Update handoff generation timestamps ---
# Overview "Latest HO Xliff Generate Date" and de-de "Latest Handoff Datetime" shared the same
# timestamp before and move to the same new timestamp together.
$wsOverview.Range("G2").Value = "2016-08-29 19:05:21"
$wsDeDe.Range("H2").Value = "2016-08-29 19:05:21"

# zh-cn "Latest Handoff Datetime" moves to its own (slightly earlier) new timestamp.
$wsZhCn.Range("H2").Value = "2016-08-29 19:05:16"

# --- Narrow the "Status" columns that used to be sized for the old, longer status text ---
$wsOverview.Columns.Item(5).ColumnWidth = 16.33   # column E (zh-cn status)
$wsOverview.Columns.Item(6).ColumnWidth = 16.33   # column F (de-de status)
$wsZhCn.Columns.Item(3).ColumnWidth = 16.33        # column C (Status)
$wsDeDe.Columns.Item(3).ColumnWidth = 16.33        # column C (Status)
